$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values ---
$ws.Range("A1").Value = "JSDSIDKDAD"

$ws.Range("E1").Value = "GWWDWDWKDKWDKKDLDLJHSJCJ"
$ws.Range("E2").Value = "SPSP[S]"
$ws.Range("E3").Value = "SPSP[S]"
$ws.Range("E4").Value = "D"
$ws.Range("E5").Value = "SD"
$ws.Range("E6").Value = "DJWJDW[DW"
$ws.Range("E7").Value = "WDJWJDND"
$ws.Range("E8").Value = "DJKDM"
$ws.Range("E9").Value = "WNNEF"
$ws.Range("E10").Value = "W"
$ws.Range("E11").Value = "MEF"
$ws.Range("E12").Value = "JWWBBWE"
$ws.Range("E13").Value = "WNBBF"
$ws.Range("E14").Value = "WJ"

# NB: order chosen to reproduce the original author's shared-string table order
$ws.Range("F7").Value = "EHDHEHEFJEE"
$ws.Range("H3").Value = "DHGBD'"
$ws.Range("H5").Value = "BVBND"
$ws.Range("G3").Value = "DBVBD"
$ws.Range("F3").Value = "DBVD"
$ws.Range("F4").Value = "DBVD"
$ws.Range("F6").Value = "DBV D"
$ws.Range("G4").Value = " BD"
$ws.Range("G5").Value = "FBVF"
$ws.Range("G6").Value = " DHHFF"
$ws.Range("G7").Value = " FJJF"
$ws.Range("H4").Value = "F   F HHJF F"
$ws.Range("H6").Value = "KKF"
$ws.Range("H7").Value = "NNNF"
